$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("D2").Value = 199
$ws1.Range("H2").Value = 16.05
$ws1.Range("L2").Value = 1.15

# Row 3
$ws1.Range("D3").Value = 228
$ws1.Range("H3").Value = 13.13
$ws1.Range("L3").Value = 1.19

# Row 4
$ws1.Range("H4").Value = 12.82
$ws1.Range("L4").Value = 1.15

# Row 5
$ws1.Range("H5").Value = 11.54
$ws1.Range("L5").Value = 1.11

# Row 6
$ws1.Range("D6").Value = 219
$ws1.Range("H6").Value = 10.64
$ws1.Range("L6").Value = 1.1

# Row 7
$ws1.Range("D7").Value = 217
$ws1.Range("H7").Value = 9.720000000000001
$ws1.Range("L7").Value = 0.96

# Row 8
$ws1.Range("H8").Value = 8.93
$ws1.Range("L8").Value = 0.88

# Row 9
$ws1.Range("H9").Value = 7.61
$ws1.Range("L9").Value = 1.05

# Row 10
$ws1.Range("H10").Value = 6.77
$ws1.Range("L10").Value = 1.07

# Row 11
$ws1.Range("H11").Value = 5.95
$ws1.Range("L11").Value = 1.12

# Row 12
$ws1.Range("D12").Value = 207
$ws1.Range("H12").Value = 5
$ws1.Range("L12").Value = 1.17

# Row 13
$ws1.Range("H13").Value = 3.93
$ws1.Range("L13").Value = 0.87

# Row 14
$ws1.Range("H14").Value = 3.03
$ws1.Range("L14").Value = 0.96

# Row 15
$ws1.Range("H15").Value = 2.16
$ws1.Range("L15").Value = 1.11

# Row 16
$ws1.Range("H16").Value = 1.15
$ws1.Range("L16").Value = 0.8

# Row 17
$ws1.Range("D17").Value = 184
$ws1.Range("H17").Value = 0.16
$ws1.Range("L17").Value = 0.96

# --- Sheet "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'3348"
$ws2.Range("B10").Value = "'1733"
$ws2.Range("B11").Value = "'864"
$ws2.Range("B12").Value = "'228"
$ws2.Range("B14").Value = "'184"
